# Update "想去人数" (want-to-go count) values in column F for a handful of
# rows across three worksheets, reflecting a refreshed data scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 2575   # was 2572
$wsExhibit.Range("F18").Value = 3689   # was 3685
$wsExhibit.Range("F24").Value = 29     # was 28
$wsExhibit.Range("F32").Value = 992    # was 991

# --- Sheet "本地生活" ----------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 261   # was 260
$wsLocal.Range("F6").Value = 35    # was 34

# --- Sheet "全部类型" ----------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value  = 261    # was 260
$wsAll.Range("F12").Value = 2575   # was 2572
$wsAll.Range("F13").Value = 2575   # was 2572
$wsAll.Range("F32").Value = 3689   # was 3685
$wsAll.Range("F37").Value = 29     # was 28
$wsAll.Range("F49").Value = 992    # was 991

$wb.Save()
